# Fill in previously-blank Manufacturer / Model / Serial No. cells for a
# handful of equipment rows, fix a manufacturer-name typo, backfill a missing
# calibration date, refresh the "duplicate values" conditional formatting on
# column E, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 - CENTRIFUGE / LABI-CE-07 / CF25073683
$ws.Cells.Item(14, 2).Value = "DLAB"
$ws.Cells.Item(14, 3).Value = "D2012 plus"
$ws.Cells.Item(14, 4).Value = "LP189AJ0000010"

# Row 16 - CENTRIFUGE / LABI-CE-09 / CF25073686
$ws.Cells.Item(16, 2).Value = "KOKUSAN"
$ws.Cells.Item(16, 3).Value = "H-19α"
$ws.Cells.Item(16, 4).Value = "153011"

# Row 26 - CENTRIFUGE / LABM-CE-03 / CF25073696
$ws.Cells.Item(26, 2).Value = "MPW"
$ws.Cells.Item(26, 3).Value = "M-DIAGNOSTIC"
$ws.Cells.Item(26, 4).Value = "102MD089221"

# Row 34 - REFRIGERATED CENTRIFUGE / LABB-RC-02 / CF25073704
$ws.Cells.Item(34, 2).Value = "ThermoFisher"
$ws.Cells.Item(34, 3).Value = "Cryofuge 6000i"
$ws.Cells.Item(34, 4).Value = "41318884"

# Row 67 - ROTATOR / LABI-RT-01 / CF25073735 (also backfill calibration date)
$ws.Cells.Item(67, 2).Value = "GEMMY"
$ws.Cells.Item(67, 3).Value = "VRN-360"
$ws.Cells.Item(67, 4).Value = "-"
$ws.Cells.Item(67, 6).Value = 45855

# Row 77 - fix manufacturer typo "LABTECH" -> "LABOTECH" (matches row 79's
# SLIDE WARMER / SW-40 combo)
$ws.Cells.Item(77, 2).Value = "LABOTECH"

# Refresh the "duplicate values" conditional formatting on column E so it
# gets its own (duplicated) differential-format record.
$rng = $ws.Range("E1:E1048576")
$rng.FormatConditions.Delete()
$fc = $rng.FormatConditions.AddUniqueValues()
$fc.DupeUnique = 1
$fc.Font.Color = 393372
$fc.Interior.Color = 13551615

# Move the active selection / scroll position.
$ws.Range("B55").Select()
